# Add IngredientRoleTypeEnum and update example with structured roles
#
# 1. IngredientRole and its subclasses (ElementSourceRole, BufferRole,
#    SolventRole, NutrientRole, VitaminRole, MineralNutrientRole) gain a
#    `role_type` column (with a dropdown data validation enumerating the
#    IngredientRoleTypeEnum permissible values) and, where missing, a
#    `source_element` column, both ahead of the trailing `owl_subclass_of`
#    column.
# 2. ProportionalPart (the RCH2 medium example) gains `composed_of` and
#    `has_role` columns up front (with `has_role` getting a dropdown data
#    validation for the ingredient role types) plus `minimal_percentage`
#    and `maximum_percentage` columns before the trailing `owl_subclass_of`
#    column.

$wb = $excel.ActiveWorkbook

$roleTypeFormula = '"carbon_source,nitrogen_source,phosphorus_source,sulfur_source,electron_donor,electron_acceptor,buffer,solvent,vitamin,mineral,trace_element,growth_factor,antibiotic,inducer,substrate"'
$ingredientRoleFormula = '"active ingredient,inactive ingredient,excipient,solvent"'

# --- IngredientRole: A1:A1 (owl_subclass_of) -> A1:C1 ---------------------
$ws = $wb.Worksheets.Item("IngredientRole")
$ws.Range("A1").Value = "source_element"
$ws.Range("B1").Value = "role_type"
$ws.Range("C1").Value = "owl_subclass_of"
$ws.Range("B2:B1048576").Validation.Add(3, 1, 1, $roleTypeFormula)

# --- ElementSourceRole: A1:B1 (source_element, owl_subclass_of) -> A1:C1 --
$ws = $wb.Worksheets.Item("ElementSourceRole")
$ws.Range("A1").Value = "source_element"
$ws.Range("B1").Value = "role_type"
$ws.Range("C1").Value = "owl_subclass_of"
$ws.Range("B2:B1048576").Validation.Add(3, 1, 1, $roleTypeFormula)

# --- BufferRole: A1:A1 (owl_subclass_of) -> A1:C1 --------------------------
$ws = $wb.Worksheets.Item("BufferRole")
$ws.Range("A1").Value = "source_element"
$ws.Range("B1").Value = "role_type"
$ws.Range("C1").Value = "owl_subclass_of"
$ws.Range("B2:B1048576").Validation.Add(3, 1, 1, $roleTypeFormula)

# --- SolventRole: A1:A1 (owl_subclass_of) -> A1:C1 -------------------------
$ws = $wb.Worksheets.Item("SolventRole")
$ws.Range("A1").Value = "source_element"
$ws.Range("B1").Value = "role_type"
$ws.Range("C1").Value = "owl_subclass_of"
$ws.Range("B2:B1048576").Validation.Add(3, 1, 1, $roleTypeFormula)

# --- NutrientRole: A1:A1 (owl_subclass_of) -> A1:C1 -------------------------
$ws = $wb.Worksheets.Item("NutrientRole")
$ws.Range("A1").Value = "source_element"
$ws.Range("B1").Value = "role_type"
$ws.Range("C1").Value = "owl_subclass_of"
$ws.Range("B2:B1048576").Validation.Add(3, 1, 1, $roleTypeFormula)

# --- VitaminRole: A1:A1 (owl_subclass_of) -> A1:C1 --------------------------
$ws = $wb.Worksheets.Item("VitaminRole")
$ws.Range("A1").Value = "source_element"
$ws.Range("B1").Value = "role_type"
$ws.Range("C1").Value = "owl_subclass_of"
$ws.Range("B2:B1048576").Validation.Add(3, 1, 1, $roleTypeFormula)

# --- MineralNutrientRole: A1:A1 (owl_subclass_of) -> A1:C1 ------------------
$ws = $wb.Worksheets.Item("MineralNutrientRole")
$ws.Range("A1").Value = "source_element"
$ws.Range("B1").Value = "role_type"
$ws.Range("C1").Value = "owl_subclass_of"
$ws.Range("B2:B1048576").Validation.Add(3, 1, 1, $roleTypeFormula)

# --- ProportionalPart: A1:C1 (has_ingredient_role, concentration, --------
#     owl_subclass_of) -> A1:G1 -------------------------------------------
$ws = $wb.Worksheets.Item("ProportionalPart")
$ws.Range("A1").Value = "composed_of"
$ws.Range("B1").Value = "has_role"
$ws.Range("C1").Value = "has_ingredient_role"
$ws.Range("D1").Value = "concentration"
$ws.Range("E1").Value = "minimal_percentage"
$ws.Range("F1").Value = "maximum_percentage"
$ws.Range("G1").Value = "owl_subclass_of"
$ws.Range("B2:B1048576").Validation.Add(3, 1, 1, $ingredientRoleFormula)

Write-Output "edits applied"
